$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relacionamento (2)")

# "aberta por" -> "aberta" on the three relationship cells that referenced it
$ws.Range("C16").Value = "aberta"
$ws.Range("C32").Value = "aberta"
$ws.Range("C46").Value = "aberta"

# New relationship row under FUNCIONARIO block (rows 31-36): FUNCIONARIO "gera" PREDITIVA IA
$ws.Range("B36").Value = "PREDITIVA IA"
$ws.Range("C36").Value = "gera"

# New relationship row under PREDITIVA IA block (rows 58-63): PREDITIVA IA "gera" FUNCIONARIO
$ws.Range("B61").Value = "FUNCIONARIO"
$ws.Range("C61").Value = "gera"

# Scroll/selection state change recorded in the diff
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("C62").Select()
